$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to append after the last existing row (row 46 -> new row 47)
$row = 47

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"

# Date column: match the number format used by the rows above (D column)
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 4).Value = 44832

$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112042
$ws.Cells.Item($row, 7).Value = "Locoto"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 80
$ws.Cells.Item($row, 11).Value = 2500
$ws.Cells.Item($row, 12).Value = 2500
$ws.Cells.Item($row, 13).Value = 2500
$ws.Cells.Item($row, 14).Value = "$/kilo"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 2500
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"

$wb.Save()
